$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The only substantive content change is that column C, which is entirely
# empty, is removed - shifting the former columns D (Habitat Type),
# E (Purpose) and F (Habitat Type 2) one position to the left (C, D, E).
$ws.Range("C1").EntireColumn.Delete()

# Update the saved selection/active cell to match the workbook's new
# sheet view state.
$ws.Range("F11").Select()
